# "update to manual status column;"
#
# The manualStatus column (I) previously stored the raw number 4 for every
# row that had been manually audited. Those cells are switched to the text
# label "[4]" instead. Column F (fastqFileName) gets a wider, explicit
# column width and the affected rows get a slightly shorter row height to
# match the new formatting, and the sheet's scroll/selection position is
# updated to the last edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose manualStatus (column I) holds the manual-audit marker.
$rows = @(3, 4, 5, 6, 7, 8, 25, 26, 27, 28, 29, 30)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "[4]"
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Give the fastqFileName column (F) room to show the full, long file names.
$ws.Columns.Item(6).ColumnWidth = 51.5

# Leave the selection on the last cell that was touched, matching the
# scrolled-down view of the edited region.
$ws.Range("I30").Select() | Out-Null
